$d = $word.ActiveDocument
$p = $d.Paragraphs(48)
$r = $p.Range
$xml = '<w:p w:rsidR="00D20962" w:rsidRPr="0033046C" w:rsidRDefault="0033046C" w:rsidP="00D20962"><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>To be industry ready over the upcoming 12 months, I have listed the following strengths that will help me achieve my goal in becoming a junior-mid programmer in the industry.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>I w</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">ill learn a different language to the point of being able to use it at industry level, whenever I have free time. </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">This will allow me to get a wider variety of jobs and not be limited by just two languages. </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Since I have experience in C# and C++, and knowledge in Unity and Unreal engine I can demonstrate an industry level prototype game to the market.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">The weakness I have can be overcome by working on a schedule to lay out a plan to work towards the bigger goal identified. Start development in C# towards a mobile app/game </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">development and add it to the professional portfolio. Use knowledge from teachers in AIE who have vast industry experience </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">to get support in creating presentable professional portfolio. Am willing to relocate to another place, if need be for a job as well as travel overseas if visa is offered by the company. Will keep uploading my projects on my portfolio and maintaining good networks with the students who are standing out more, so I can be recommended by them whilst I compete with them. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>As support I can use the opportunities available to me, such as VR pixels which is an event that is promoted to students in AIE. This will allow me to increase knowledge on games and increase contacts and good relations with people from the industry whilst increasing my social skills. To increase my competitiveness, I will keep in contact with my teachers, and keep asking for good criticism upon my portfolio to make it more presentable and better than others. Have multiple books on physics from high school and still have contact with my physics teachers, which helps in physics-based calculations that are difficult.</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> To</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> identify, the</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> potential threats</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> are competition, visa, travel, time and networking ability. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>I will minimise the impact of these threats by preparing myself for them in advance. I will take into consideration of the visa, and travel duratio</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">n and time to relocate to the other place before applying for the job </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>and</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> focus my search on a larger more stable company in the capital cities I have a visa to work in. I can minimise the impact of competitiveness by joining and starting off with entry level jobs, and focus on my </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>strengths</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> such as C++, C# and development in Unity in the Sydney area. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">In conclusion, these will allow me to maximize my ability to obtain the selected job and be prepared industry jobs in 12 months’ time.  </w:t></w:r></w:p>'
$r.InsertXML($xml)
